# The deck's two theme parts (ppt/theme/theme1.xml, used by the slide
# master, and ppt/theme/theme2.xml, used by the notes master) had their
# contents swapped: theme1 ("Integral") <-> theme2 ("Office Theme").
#
# Re-apply that swap by writing the target ("Office Theme") colour
# scheme values into the presentation's live theme colour scheme
# (ppt/theme/theme1.xml's <a:clrScheme>), via the slide's
# ThemeColorScheme, which is the color-scheme surface the object model
# exposes for editing theme colours.
#
# Index order (1-based) of ThemeColorScheme.Item(n) matches the OOXML
# <a:clrScheme> child order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      -> #000000
$tcs.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      -> #FFFFFF
$tcs.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      -> #44546A
$tcs.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      -> #E7E6E6
$tcs.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  -> #5B9BD5
$tcs.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  -> #ED7D31
$tcs.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  -> #A5A5A5
$tcs.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  -> #FFC000
$tcs.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  -> #4472C4
$tcs.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  -> #70AD47
$tcs.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    -> #0563C1
$tcs.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink -> #954F72
